# Final submission update for organizacion.xlsx
# - Rename the "Fecha" column header to "Fecha fin e inicio"
# - Add a new task row for Gaizka: "Mostrar las herramientas por tarea (many to many)"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header of column D (was "Fecha")
$ws.Range("D1").Value = "Fecha fin e inicio"

# Append the new row (row 31) with the same layout/format as the previous rows
$ws.Range("B31").Value = "Gaizka"
$ws.Range("C31").Value = "Mostrar las herramientas por tarea (many to many)"
$ws.Range("D31").Value = 45784
$ws.Range("D31").NumberFormat = $ws.Range("D30").NumberFormat

# Leave the cursor where the author would have ended up after typing the new row
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$ws.Range("B32").Select()
